$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# Locate the sentence that is being edited. Find.Execute() (with no
# replacement text) collapses/extends the range it is called on to the
# exact span of the match, so $target ends up bounding precisely
# "Mon choix\u00A0: Novembre." (paragraph mark excluded).
$target = $d.Content
$found = $target.Find.Execute("Mon choix : Novembre.", $true, $false, $false,
                               $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the sentence 'Mon choix : Novembre.' to edit."
}

# Replace that span with the same sentence split across three runs, and
# with "Novembre" lower-cased to "novembre" - i.e. exactly the OOXML
# produced by the diff:
#   <w:r><w:t>Mon choix :</w:t></w:r>
#   <w:r><w:t xml:space="preserve"> n</w:t></w:r>
#   <w:r><w:t>ovembre.</w:t></w:r>
# InsertXML replaces the contents of the target range with the supplied
# WordprocessingML, which is how Word COM creates/splits runs from a
# script (ordinary .Text assignment would just rewrite the text of the
# existing single run).
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>Mon choix${nbsp}:</w:t></w:r>
            <w:r><w:t xml:space="preserve"> n</w:t></w:r>
            <w:r><w:t>ovembre.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$target.InsertXML($xml)
